# Updates cryptos list price (D) and hourly volume change (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.610.16'
$ws.Range('E2').Value = '  -2.60%  '
$ws.Range('D3').Value = '3.923.40'
$ws.Range('E3').Value = '  -2.57%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '540.37'
$ws.Range('E5').Value = '  +5.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.73'
$ws.Range('E6').Value = '  +0.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.685'
$ws.Range('E7').Value = '  -6.00%  '
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.734'
$ws.Range('E9').Value = '  -4.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.167'
$ws.Range('E10').Value = '  -3.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.95'
$ws.Range('E11').Value = '  +14.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000317'
$ws.Range('E12').Value = '  -1.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.47'
$ws.Range('E13').Value = '  -2.52%  '
$ws.Range('D14').Value = '4.564.18'
$ws.Range('E14').Value = '  -2.08%  '
$ws.Range('D15').Value = '3.929.15'
$ws.Range('E15').Value = '  -1.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.87'
$ws.Range('E16').Value = '  -0.93%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.18'
$ws.Range('E17').Value = '  -3.78%  '
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.17'
$ws.Range('E19').Value = '  -2.79%  '
$ws.Range('D20').Value = '70.547.48'
$ws.Range('E20').Value = '  -2.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '424.50'
$ws.Range('E21').Value = '  -2.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '95.80'
$ws.Range('E22').Value = '  -6.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.50'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.22'
$ws.Range('E24').Value = '  +7.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.17'
$ws.Range('E25').Value = '  -2.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.17'
$ws.Range('E26').Value = '  -2.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.45'
$ws.Range('E27').Value = '  -4.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.84'
$ws.Range('E28').Value = '  +0.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.60'
$ws.Range('E29').Value = '  +15.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.14'
$ws.Range('E30').Value = '  -3.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.36'
$ws.Range('E31').Value = '  +8.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '13.31'
$ws.Range('E32').Value = '  -1.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.128'
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '676.46'
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '47.45'
$ws.Range('E35').Value = '  +16.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '65.08'
$ws.Range('E36').Value = '  -3.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.430'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').Value = '0.0₃0814'
$ws.Range('E38').Value = '  -5.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.148'
$ws.Range('E39').Value = '  -1.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.35'
$ws.Range('E40').Value = '  -4.57%  '
$ws.Range('E41').Value = '  +5.30%  '
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0478'
$ws.Range('E44').Value = '  -1.00%  '
$ws.Range('E45').Value = '  -6.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.68'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.57'
$ws.Range('E47').Value = '  +6.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.36'
$ws.Range('E48').Value = '  -3.95%  '
$ws.Range('E49').Value = '  -3.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000272'
$ws.Range('E50').Value = '  +2.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '144.49'
$ws.Range('E51').Value = '  +1.67%  '
